# Applies the crypto price/volume/coin-listing refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Val = '26.823.88' },
    @{ Row = 2; Col = 5; Val = '  +3.96%  ' },
    @{ Row = 3; Col = 4; Val = '1.878.22' },
    @{ Row = 3; Col = 5; Val = '  +3.43%  ' },
    @{ Row = 4; Col = 5; Val = '  +0.07%  ' },
    @{ Row = 5; Col = 4; Val = '277.24' },
    @{ Row = 5; Col = 5; Val = '  +0.16%  ' },
    @{ Row = 6; Col = 4; Val = '1.001' },
    @{ Row = 6; Col = 5; Val = '  +0.04%  ' },
    @{ Row = 7; Col = 4; Val = '0.5287' },
    @{ Row = 7; Col = 5; Val = '  +3.78%  ' },
    @{ Row = 8; Col = 4; Val = '0.3425' },
    @{ Row = 8; Col = 5; Val = '  -2.93%  ' },
    @{ Row = 9; Col = 2; Val = 'Dogecoin' },
    @{ Row = 9; Col = 3; Val = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge' },
    @{ Row = 9; Col = 4; Val = '0.06945' },
    @{ Row = 9; Col = 5; Val = '  +4.16%  ' },
    @{ Row = 10; Col = 2; Val = 'Solana' },
    @{ Row = 10; Col = 3; Val = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol' },
    @{ Row = 10; Col = 4; Val = '20.01' },
    @{ Row = 10; Col = 5; Val = '  -0.04%  ' },
    @{ Row = 11; Col = 2; Val = 'Polygon' },
    @{ Row = 11; Col = 3; Val = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' },
    @{ Row = 11; Col = 4; Val = '0.8025' },
    @{ Row = 11; Col = 5; Val = '  -3.20%  ' },
    @{ Row = 12; Col = 2; Val = 'TRON' },
    @{ Row = 12; Col = 3; Val = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx' },
    @{ Row = 12; Col = 4; Val = '0.07740' },
    @{ Row = 12; Col = 5; Val = '  -1.51%  ' },
    @{ Row = 13; Col = 2; Val = 'WrappedEther' },
    @{ Row = 13; Col = 3; Val = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' },
    @{ Row = 13; Col = 4; Val = '1.889.92' },
    @{ Row = 13; Col = 5; Val = '  +4.12%  ' },
    @{ Row = 14; Col = 2; Val = 'Litecoin' },
    @{ Row = 14; Col = 3; Val = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc' },
    @{ Row = 14; Col = 4; Val = '90.21' },
    @{ Row = 14; Col = 5; Val = '  +3.12%  ' },
    @{ Row = 15; Col = 2; Val = 'Polkadot' },
    @{ Row = 15; Col = 3; Val = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Row = 15; Col = 4; Val = '5.168' },
    @{ Row = 15; Col = 5; Val = '  +1.74%  ' },
    @{ Row = 16; Col = 2; Val = 'Avalanche' },
    @{ Row = 16; Col = 3; Val = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' },
    @{ Row = 16; Col = 4; Val = '14.54' },
    @{ Row = 16; Col = 5; Val = '  +2.90%  ' },
    @{ Row = 17; Col = 2; Val = 'BinanceUSD' },
    @{ Row = 17; Col = 3; Val = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' },
    @{ Row = 17; Col = 4; Val = '1.001' },
    @{ Row = 17; Col = 5; Val = '  +0.06%  ' },
    @{ Row = 18; Col = 2; Val = 'ShibaInu' },
    @{ Row = 18; Col = 3; Val = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib' },
    @{ Row = 18; Col = 4; Val = '0.000008052' },
    @{ Row = 18; Col = 5; Val = '  +0.10%  ' },
    @{ Row = 19; Col = 2; Val = 'Dai' },
    @{ Row = 19; Col = 3; Val = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' },
    @{ Row = 19; Col = 4; Val = '1.001' },
    @{ Row = 19; Col = 5; Val = '  +0.05%  ' },
    @{ Row = 20; Col = 2; Val = 'WrappedBTC' },
    @{ Row = 20; Col = 3; Val = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Row = 20; Col = 4; Val = '26.876.41' },
    @{ Row = 20; Col = 5; Val = '  +3.97%  ' },
    @{ Row = 21; Col = 2; Val = 'WrappedliquidstakedEther2.0' },
    @{ Row = 21; Col = 3; Val = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' },
    @{ Row = 21; Col = 4; Val = '2.129.02' },
    @{ Row = 21; Col = 5; Val = '  +4.03%  ' },
    @{ Row = 22; Col = 2; Val = 'Uniswap' },
    @{ Row = 22; Col = 3; Val = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' },
    @{ Row = 22; Col = 4; Val = '4.736' },
    @{ Row = 22; Col = 5; Val = '  +0.16%  ' },
    @{ Row = 23; Col = 2; Val = 'Cosmos' },
    @{ Row = 23; Col = 3; Val = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' },
    @{ Row = 23; Col = 4; Val = '10.01' },
    @{ Row = 23; Col = 5; Val = '  +0.01%  ' },
    @{ Row = 24; Col = 2; Val = 'Chainlink' },
    @{ Row = 24; Col = 3; Val = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link' },
    @{ Row = 24; Col = 4; Val = '6.195' },
    @{ Row = 24; Col = 5; Val = '  +1.74%  ' },
    @{ Row = 25; Col = 2; Val = 'LidoDAOToken' },
    @{ Row = 25; Col = 3; Val = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo' },
    @{ Row = 25; Col = 4; Val = '2.377' },
    @{ Row = 25; Col = 5; Val = '  +8.27%  ' },
    @{ Row = 26; Col = 2; Val = 'Monero' },
    @{ Row = 26; Col = 3; Val = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' },
    @{ Row = 26; Col = 4; Val = '146.59' },
    @{ Row = 26; Col = 5; Val = '  +3.82%  ' },
    @{ Row = 27; Col = 2; Val = 'Toncoin' },
    @{ Row = 27; Col = 3; Val = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton' },
    @{ Row = 27; Col = 4; Val = '1.667' },
    @{ Row = 27; Col = 5; Val = '  -0.33%  ' },
    @{ Row = 28; Col = 2; Val = 'EthereumClassic' },
    @{ Row = 28; Col = 3; Val = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' },
    @{ Row = 28; Col = 4; Val = '17.33' },
    @{ Row = 28; Col = 5; Val = '  +1.50%  ' },
    @{ Row = 29; Col = 2; Val = 'BitcoinCash' },
    @{ Row = 29; Col = 3; Val = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch' },
    @{ Row = 29; Col = 4; Val = '113.53' },
    @{ Row = 29; Col = 5; Val = '  +3.82%  ' },
    @{ Row = 30; Col = 2; Val = 'InternetComputer(DFINITY)' },
    @{ Row = 30; Col = 3; Val = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' },
    @{ Row = 30; Col = 4; Val = '4.338' },
    @{ Row = 30; Col = 5; Val = '  -0.33%  ' },
    @{ Row = 31; Col = 2; Val = 'Filecoin' },
    @{ Row = 31; Col = 3; Val = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil' },
    @{ Row = 31; Col = 4; Val = '4.295' },
    @{ Row = 31; Col = 5; Val = '  +1.27%  ' },
    @{ Row = 32; Col = 2; Val = 'Stellar' },
    @{ Row = 32; Col = 3; Val = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm' },
    @{ Row = 32; Col = 4; Val = '0.08885' },
    @{ Row = 32; Col = 5; Val = '  +0.99%  ' },
    @{ Row = 33; Col = 2; Val = 'Hedera' },
    @{ Row = 33; Col = 3; Val = 'https://coinranking.com/coin/jad286TjB+hedera-hbar' },
    @{ Row = 33; Col = 4; Val = '0.04891' },
    @{ Row = 33; Col = 5; Val = '  -0.18%  ' },
    @{ Row = 34; Col = 2; Val = 'ARBITRUM' },
    @{ Row = 34; Col = 3; Val = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb' },
    @{ Row = 34; Col = 4; Val = '1.174' },
    @{ Row = 34; Col = 5; Val = '  +3.26%  ' },
    @{ Row = 35; Col = 2; Val = 'ImmutableX' },
    @{ Row = 35; Col = 3; Val = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' },
    @{ Row = 35; Col = 4; Val = '0.7248' },
    @{ Row = 35; Col = 5; Val = '  -0.42%  ' },
    @{ Row = 36; Col = 2; Val = 'HuobiToken' },
    @{ Row = 36; Col = 3; Val = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' },
    @{ Row = 36; Col = 4; Val = '2.892' },
    @{ Row = 36; Col = 5; Val = '  +0.95%  ' },
    @{ Row = 37; Col = 2; Val = 'MXToken' },
    @{ Row = 37; Col = 3; Val = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' },
    @{ Row = 37; Col = 4; Val = '3.287' },
    @{ Row = 37; Col = 5; Val = '  +4.93%  ' },
    @{ Row = 38; Col = 2; Val = 'RenderToken' },
    @{ Row = 38; Col = 3; Val = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' },
    @{ Row = 38; Col = 4; Val = '2.344' },
    @{ Row = 38; Col = 5; Val = '  -1.25%  ' },
    @{ Row = 39; Col = 2; Val = 'VeChain' },
    @{ Row = 39; Col = 3; Val = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Row = 39; Col = 4; Val = '0.01838' },
    @{ Row = 39; Col = 5; Val = '  -0.74%  ' },
    @{ Row = 40; Col = 2; Val = 'TheSandbox' },
    @{ Row = 40; Col = 3; Val = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand' },
    @{ Row = 40; Col = 4; Val = '0.5103' },
    @{ Row = 40; Col = 5; Val = '  -1.67%  ' },
    @{ Row = 41; Col = 2; Val = 'TrustWalletToken' },
    @{ Row = 41; Col = 3; Val = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' },
    @{ Row = 41; Col = 4; Val = '0.9550' },
    @{ Row = 41; Col = 5; Val = '  -0.20%  ' },
    @{ Row = 42; Col = 2; Val = 'Quant' },
    @{ Row = 42; Col = 3; Val = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt' },
    @{ Row = 42; Col = 4; Val = '116.07' },
    @{ Row = 42; Col = 5; Val = '  +4.47%  ' },
    @{ Row = 43; Col = 2; Val = 'FraxShare' },
    @{ Row = 43; Col = 3; Val = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' },
    @{ Row = 43; Col = 4; Val = '6.165' },
    @{ Row = 43; Col = 5; Val = '  -0.64%  ' },
    @{ Row = 44; Col = 2; Val = 'Aptos' },
    @{ Row = 44; Col = 3; Val = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' },
    @{ Row = 44; Col = 4; Val = '8.064' },
    @{ Row = 44; Col = 5; Val = '  +0.76%  ' },
    @{ Row = 45; Col = 2; Val = 'PaxDollar' },
    @{ Row = 45; Col = 3; Val = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp' },
    @{ Row = 45; Col = 4; Val = '1.000' },
    @{ Row = 45; Col = 5; Val = '  +0.01%  ' },
    @{ Row = 46; Col = 2; Val = 'Decentraland' },
    @{ Row = 46; Col = 3; Val = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' },
    @{ Row = 46; Col = 4; Val = '0.4447' },
    @{ Row = 46; Col = 5; Val = '  -2.70%  ' },
    @{ Row = 47; Col = 2; Val = 'Algorand' },
    @{ Row = 47; Col = 3; Val = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo' },
    @{ Row = 47; Col = 4; Val = '0.1338' },
    @{ Row = 47; Col = 5; Val = '  -2.06%  ' },
    @{ Row = 48; Col = 2; Val = 'EnergySwap' },
    @{ Row = 48; Col = 3; Val = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' },
    @{ Row = 48; Col = 4; Val = '9.303' },
    @{ Row = 48; Col = 5; Val = '  +0.71%  ' },
    @{ Row = 49; Col = 2; Val = 'Elrond' },
    @{ Row = 49; Col = 3; Val = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' },
    @{ Row = 49; Col = 4; Val = '36.12' },
    @{ Row = 49; Col = 5; Val = '  -1.39%  ' },
    @{ Row = 50; Col = 2; Val = 'Cronos' },
    @{ Row = 50; Col = 3; Val = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro' },
    @{ Row = 50; Col = 4; Val = '0.05945' },
    @{ Row = 50; Col = 5; Val = '  +1.82%  ' },
    @{ Row = 51; Col = 2; Val = 'NEARProtocol' },
    @{ Row = 51; Col = 3; Val = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near' },
    @{ Row = 51; Col = 4; Val = '1.489' },
    @{ Row = 51; Col = 5; Val = '  -0.91%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    # Force text format so numeric-looking strings (e.g. "1.001", "0.07740")
    # are preserved exactly as text instead of being parsed into numbers.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Val
    # Reset to the default style so we do not leave a stray text-format style
    # applied to the cell (matches original formatting of these cells).
    $cell.Style = "Normal"
}

Write-Host "Applied $($updates.Count) cell updates."
